# Fix Training Data Issue (#48)
# The "Date" column (BF) held a malformed label ("6-21-2012-13") resulting
# from the way NBA stats were scraped — one day off. Correct it to an
# ISO-8601 date string ("2013-06-21") for every data row (BF2:BF31).
# BF1 is the "Date" header and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("BF2:BF31")
# Force text storage so Excel doesn't reinterpret the ISO date string as a
# date serial number when we assign it below.
$dataRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # Column BF = 58
    $cell.Value = "2013-06-21"
}
